$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.668431282043457
$ws.Range("B1").Value = 4.453009128570557
$ws.Range("C1").Value = 3.409494638442993
$ws.Range("D1").Value = 2.592077016830444
$ws.Range("E1").Value = 2.167780160903931
